$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row reorderings (coin rank swaps): update Coin name, Link, Price, Volume ---
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''2.297'
$ws.Range("E46").Value = '  +2.88%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.6361'
$ws.Range("E47").Value = '  -1.49%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '''1.245'
$ws.Range("E49").Value = '  -2.43%  '

$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '''1.227'
$ws.Range("E50").Value = '  +6.13%  '

# --- Price / Volume(1h) updates ---
$ws.Range("D2").Value = '30.501.01'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '2.102.31'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '''332.27'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '''0.5227'
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("D8").Value = '''0.4485'
$ws.Range("E8").Value = '  +2.34%  '
$ws.Range("D9").Value = '''53.52'
$ws.Range("E9").Value = '  +16.66%  '
$ws.Range("D10").Value = '''0.08923'
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").Value = '''1.155'
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("D12").Value = '''24.44'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '2.094.34'
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("D14").Value = '''6.728'
$ws.Range("D15").Value = '''7.709'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").Value = '''96.33'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").Value = '''1.004'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = '''0.00001124'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").Value = '''0.06623'
$ws.Range("E19").Value = '  -0.73%  '
$ws.Range("D20").Value = '''19.21'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '''6.287'
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("D23").Value = '30.560.15'
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '''12.33'
$ws.Range("E24").Value = '  +1.49%  '
$ws.Range("D25").Value = '''2.322'
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").Value = '2.342.79'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '''22.29'
$ws.Range("E27").Value = '  -2.18%  '
$ws.Range("D28").Value = '''2.580'
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").Value = '''163.83'
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '''132.28'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").Value = '''1.197'
$ws.Range("E31").Value = '  +2.33%  '
$ws.Range("D32").Value = '''0.1073'
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").Value = '''1.672'
$ws.Range("E33").Value = '  +8.57%  '
$ws.Range("D34").Value = '''6.155'
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").Value = '''3.903'
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("E36").Value = '  +10.04%  '
$ws.Range("D37").Value = '''0.02572'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").Value = '''0.06768'
$ws.Range("E38").Value = '  +0.38%  '
$ws.Range("D39").Value = '''5.481'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").Value = '''0.2263'
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("D42").Value = '''0.6923'
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("D43").Value = '''1.256'
$ws.Range("E43").Value = '  +0.55%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''13.98'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D48").Value = '''3.639'
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D51").Value = '''82.05'
$ws.Range("E51").Value = '  -0.41%  '
